$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R values for rows 4-34 (2021 data), mirroring the style of column Q in the same row
$values = @{
    4  = 2021
    5  = 11.9
    6  = 13.1
    7  = 10.6
    8  = 11
    9  = 10
    10 = 12
    11 = 10.199999999999999
    12 = 10.5
    13 = 10
    14 = 19.399999999999999
    15 = 22.3
    16 = 16.399999999999999
    17 = 9.4
    18 = 11.4
    19 = 7.3
    20 = 3.1
    21 = 2.9
    22 = 3.4
    23 = 15
    24 = 17.3
    25 = 12.7
    26 = 7.9
    27 = 8.4
    28 = 7.4
    29 = 15.2
    30 = 17.600000000000001
    31 = 12.6
    32 = 27.9
    33 = 32.700000000000003
    34 = 22.8
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Cells.Item($row, 17)   # Q column
    $dstCell = $ws.Cells.Item($row, 18)   # R column
    $dstCell.Value = $values[$row]
    $dstCell.NumberFormat = $srcCell.NumberFormat
    $dstCell.Font.Name = $srcCell.Font.Name
    $dstCell.Font.Size = $srcCell.Font.Size
    $dstCell.Font.Bold = $srcCell.Font.Bold
    $dstCell.Borders.Item(9).LineStyle = $srcCell.Borders.Item(9).LineStyle
    $dstCell.Borders.Item(9).Weight = $srcCell.Borders.Item(9).Weight
}

$ws.Range("R3").Select()
